# Insert a new weekly price record as row 100 in the "Poroto verde" price
# history sheet. This shifts the existing rows 100-177 down to 101-178
# (extending the table to A1:R178), and populates the newly inserted row
# with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 100, pushing rows 100..177 down to 101..178
$ws.Rows("100").Insert()

# Populate the new row 100 with the new record
$ws.Cells.Item(100, 1).Value = 4
$ws.Cells.Item(100, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(100, 3).Value = "Los Lagos"
$ws.Cells.Item(100, 4).Value = 45216
$ws.Cells.Item(100, 5).Value = 10
$ws.Cells.Item(100, 6).Value = 100112031
$ws.Cells.Item(100, 7).Value = "Poroto verde"
$ws.Cells.Item(100, 8).Value = "Magnum"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 45
$ws.Cells.Item(100, 11).Value = 37000
$ws.Cells.Item(100, 12).Value = 37000
$ws.Cells.Item(100, 13).Value = 37000
$ws.Cells.Item(100, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(100, 15).Value = "Perú"
$ws.Cells.Item(100, 16).Value = 1480
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = "Hortaliza"
